$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I5").Value = 4.33
$ws.Range("J5").Value = 2.88
$ws.Range("L5").Value = 5.5
$ws.Range("M5").Value = 1.14
$ws.Range("N5").Value = 5.5
$ws.Range("AA5").Value = 4.75
$ws.Range("AH5").Value = 6.5
$ws.Range("H6").Value = 3.1
$ws.Range("I6").Value = 3.9
$ws.Range("M6").Value = 1.13
$ws.Range("N6").Value = 6
$ws.Range("Z8").Value = 1.5
$ws.Range("M9").Value = 1.08
$ws.Range("N9").Value = 7.5
$ws.Range("U9").Value = 4.33
$ws.Range("V9").Value = 1.2
$ws.Range("Z9").Value = 1.5
$ws.Range("Z10").Value = 1.67
$ws.Range("G11").Value = 2.85
$ws.Range("H11").Value = 2.65
$ws.Range("I11").Value = 2.85
$ws.Range("J11").Value = 3.55
$ws.Range("K11").Value = 1.83
$ws.Range("L11").Value = 3.5
$ws.Range("M11").Value = 1.14
$ws.Range("N11").Value = 5
$ws.Range("O11").Value = 1.55
$ws.Range("P11").Value = 2.3
$ws.Range("Q11").Value = 2.62
$ws.Range("R11").Value = 1.44
$ws.Range("U11").Value = 4.7
$ws.Range("V11").Value = 1.15
$ws.Range("W11").Value = 1.57
$ws.Range("X11").Value = 2.25
$ws.Range("Y11").Value = 2.05
$ws.Range("Z11").Value = 1.7
$ws.Range("AA11").Value = 6.5
$ws.Range("AC11").Value = 10.75
$ws.Range("AF11").Value = 50
$ws.Range("AG11").Value = 5
$ws.Range("AI11").Value = 16.5
$ws.Range("AJ11").Value = 110
$ws.Range("AL11").Value = 6.6
$ws.Range("AN11").Value = 10.75
$ws.Range("AP11").Value = 30
$ws.Range("AQ11").Value = 45
$ws.Range("G12").Value = 3.2
$ws.Range("I12").Value = 2.5
$ws.Range("J12").Value = 3.85
$ws.Range("Q12").Value = 2.5
$ws.Range("R12").Value = 1.47
$ws.Range("X12").Value = 2.32
$ws.Range("AA12").Value = 7.2
$ws.Range("AC12").Value = 11.5
$ws.Range("AM12").Value = 11.5
$ws.Range("AN12").Value = 9.5
$ws.Range("AO12").Value = 28
$ws.Range("AP12").Value = 24
$ws.Range("G14").Value = 1.35
$ws.Range("H14").Value = 4.3
$ws.Range("I14").Value = 9.5
$ws.Range("J14").Value = 1.82
$ws.Range("K14").Value = 2.27
$ws.Range("L14").Value = 8.25
$ws.Range("M14").Value = 1.07
$ws.Range("N14").Value = 7
$ws.Range("O14").Value = 1.33
$ws.Range("P14").Value = 3.05
$ws.Range("Q14").Value = 1.95
$ws.Range("R14").Value = 1.75
$ws.Range("U14").Value = 3.3
$ws.Range("V14").Value = 1.29
$ws.Range("W14").Value = 1.4
$ws.Range("X14").Value = 2.75
$ws.Range("Y14").Value = 2.37
$ws.Range("Z14").Value = 1.52
$ws.Range("AA14").Value = 5.3
$ws.Range("AB14").Value = 5.4
$ws.Range("AC14").Value = 8.75
$ws.Range("AD14").Value = 7.8
$ws.Range("AE14").Value = 13
$ws.Range("AG14").Value = 7
$ws.Range("AH14").Value = 8.75
$ws.Range("AI14").Value = 28
$ws.Range("AJ14").Value = 175
$ws.Range("AL14").Value = 18.5
$ws.Range("AM14").Value = 65
$ws.Range("AN14").Value = 30
$ws.Range("AO14").Value = 300
$ws.Range("AP14").Value = 150
$ws.Range("Q15").Value = 1.98
$ws.Range("R15").Value = 1.88
$ws.Range("U15").Value = 3.4
$ws.Range("V15").Value = 1.33
$ws.Range("I16").Value = 6
$ws.Range("J16").Value = 2.25
$ws.Range("M16").Value = 1.07
$ws.Range("N16").Value = 9
$ws.Range("W16").Value = 1.5
$ws.Range("X16").Value = 2.5
$ws.Range("Y16").Value = 2.2
$ws.Range("Z16").Value = 1.62
$ws.Range("AB16").Value = 6.5
$ws.Range("AC16").Value = 9
$ws.Range("AD16").Value = 11
$ws.Range("AJ16").Value = 81
$ws.Range("H17").Value = 3
$ws.Range("I17").Value = 2.25
$ws.Range("K17").Value = 1.91
$ws.Range("L17").Value = 3.1
$ws.Range("W17").Value = 1.57
$ws.Range("X17").Value = 2.25
$ws.Range("Y17").Value = 2.1
$ws.Range("Z17").Value = 1.67
$ws.Range("AG17").Value = 6.5
$ws.Range("AI17").Value = 19
$ws.Range("AL17").Value = 6
$ws.Range("AN17").Value = 10
$ws.Range("AQ17").Value = 41
$ws.Range("AR17").Value = 1.93
$ws.Range("AS17").Value = 1.93
$ws.Range("K20").Value = 2.1
$ws.Range("Q20").Value = 2
$ws.Range("R20").Value = 1.85
$ws.Range("U20").Value = 3.5
$ws.Range("V20").Value = 1.3
$ws.Range("AA20").Value = 8
$ws.Range("AE20").Value = 19
$ws.Range("AF20").Value = 29
$ws.Range("AG20").Value = 10
$ws.Range("AI20").Value = 15
$ws.Range("AJ20").Value = 51
$ws.Range("AK20").Value = 251
$ws.Range("AL20").Value = 9.5
$ws.Range("AM20").Value = 15
$ws.Range("AP20").Value = 26
$ws.Range("AR23").Value = 1.9
$ws.Range("AS23").Value = 1.9
$ws.Range("Q26").Value = 2.1
$ws.Range("R26").Value = 1.7
